# Auto-generated Excel COM-interop script
# Applies updated market-price / leve-profit values to several sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 884.25
$ws.Range("I4").Value = 750.5
$ws.Range("J4").Value = 1018
$ws.Range("K4").Value = 750.5
$ws.Range("L4").Value = 1018
$ws.Range("M4").Value = -636.5
$ws.Range("N4").Value = -1246
$ws.Range("H18").Value = 290.8
$ws.Range("I18").Value = 290.8
$ws.Range("K18").Value = 290.8
$ws.Range("M18").Value = -6.800000000000011
$ws.Range("H19").Value = 680.7
$ws.Range("I19").Value = 646.1539
$ws.Range("J19").Value = 707.1177
$ws.Range("K19").Value = 646.1539
$ws.Range("L19").Value = 707.1177
$ws.Range("M19").Value = -471.1539
$ws.Range("N19").Value = -1057.1177
$ws.Range("H32").Value = 11364648
$ws.Range("I32").Value = 25000386
$ws.Range("J32").Value = 1533.3334
$ws.Range("K32").Value = 25000386
$ws.Range("L32").Value = 1533.3334
$ws.Range("M32").Value = -25000060
$ws.Range("N32").Value = -2185.3334
$ws.Range("H41").Value = 278.57144
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 337.5
$ws.Range("K41").Value = 200
$ws.Range("L41").Value = 337.5
$ws.Range("M41").Value = 240
$ws.Range("N41").Value = -1217.5
$ws.Range("H51").Value = 19913.818
$ws.Range("I51").Value = 100001
$ws.Range("J51").Value = 2116.6667
$ws.Range("K51").Value = 100001
$ws.Range("L51").Value = 2116.6667
$ws.Range("M51").Value = -99517
$ws.Range("N51").Value = -3084.6667
$ws.Range("H53").Value = 22727980
$ws.Range("I53").Value = 41666924
$ws.Range("J53").Value = 1247.6
$ws.Range("K53").Value = 41666924
$ws.Range("L53").Value = 1247.6
$ws.Range("M53").Value = -41666287
$ws.Range("N53").Value = -2521.6
$ws.Range("H55").Value = 868419.5600000001
$ws.Range("I55").Value = 591.0526
$ws.Range("K55").Value = 591.0526
$ws.Range("M55").Value = -377.0526
$ws.Range("H88").Value = 11410535
$ws.Range("J88").Value = 13040269
$ws.Range("L88").Value = 13040269
$ws.Range("N88").Value = -13041081
$ws.Range("H91").Value = 11410535
$ws.Range("J91").Value = 13040269
$ws.Range("L91").Value = 13040269
$ws.Range("N91").Value = -13043077
$ws.Range("H98").Value = 96155096
$ws.Range("I98").Value = 113637680
$ws.Range("K98").Value = 113637680
$ws.Range("M98").Value = -113636182
$ws.Range("H116").Value = 2330.862
$ws.Range("I116").Value = 2090.9524
$ws.Range("J116").Value = 2960.625
$ws.Range("K116").Value = 2090.9524
$ws.Range("L116").Value = 2960.625
$ws.Range("M116").Value = 1351.0476
$ws.Range("N116").Value = -9844.625
$ws.Range("H122").Value = 96155096
$ws.Range("I122").Value = 113637680
$ws.Range("K122").Value = 340913040
$ws.Range("M122").Value = -340910590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4623.913
$ws.Range("I7").Value = 211.38461
$ws.Range("J7").Value = 10360.2
$ws.Range("K7").Value = 211.38461
$ws.Range("L7").Value = 10360.2
$ws.Range("M7").Value = -98.38461000000001
$ws.Range("N7").Value = -10586.2
$ws.Range("H86").Value = 32634312
$ws.Range("I86").Value = 55577464
$ws.Range("J86").Value = 30885.105
$ws.Range("K86").Value = 55577464
$ws.Range("L86").Value = 30885.105
$ws.Range("M86").Value = -55576341
$ws.Range("N86").Value = -33131.105
$ws.Range("H89").Value = 32634312
$ws.Range("I89").Value = 55577464
$ws.Range("J89").Value = 30885.105
$ws.Range("K89").Value = 277887320
$ws.Range("L89").Value = 154425.525
$ws.Range("M89").Value = -277881704
$ws.Range("N89").Value = -165657.525

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 400
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -3338
$ws.Range("H22").Value = 4006799
$ws.Range("I22").Value = 100000000
$ws.Range("J22").Value = 7082.25
$ws.Range("K22").Value = 300000000
$ws.Range("L22").Value = 21246.75
$ws.Range("M22").Value = -299999831
$ws.Range("N22").Value = -21584.75
$ws.Range("H27").Value = 4006799
$ws.Range("I27").Value = 100000000
$ws.Range("J27").Value = 7082.25
$ws.Range("K27").Value = 300000000
$ws.Range("L27").Value = 21246.75
$ws.Range("M27").Value = -299999898
$ws.Range("N27").Value = -21450.75
$ws.Range("H68").Value = 7289.933
$ws.Range("I68").Value = 449.875
$ws.Range("J68").Value = 15107.143
$ws.Range("K68").Value = 1349.625
$ws.Range("L68").Value = 45321.429
$ws.Range("M68").Value = -538.625
$ws.Range("N68").Value = -46943.429
$ws.Range("H69").Value = 10678
$ws.Range("J69").Value = 11603.272
$ws.Range("L69").Value = 34809.81600000001
$ws.Range("N69").Value = -36431.81600000001
$ws.Range("H71").Value = 7289.933
$ws.Range("I71").Value = 449.875
$ws.Range("J71").Value = 15107.143
$ws.Range("K71").Value = 4048.875
$ws.Range("L71").Value = 135964.287
$ws.Range("M71").Value = 7.125
$ws.Range("N71").Value = -144076.287
$ws.Range("H72").Value = 10678
$ws.Range("J72").Value = 11603.272
$ws.Range("L72").Value = 104429.448
$ws.Range("N72").Value = -112541.448
$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H86").Value = 260
$ws.Range("I86").Value = 260
$ws.Range("K86").Value = 780
$ws.Range("M86").Value = 406
$ws.Range("H89").Value = 260
$ws.Range("I89").Value = 260
$ws.Range("K89").Value = 2340
$ws.Range("M89").Value = 3588
$ws.Range("H97").Value = 3243.625
$ws.Range("J97").Value = 2992.7144
$ws.Range("L97").Value = 8978.143199999999
$ws.Range("N97").Value = -9970.143199999999
$ws.Range("H120").Value = 6716.25
$ws.Range("I120").Value = 3432.5
$ws.Range("J120").Value = 10000
$ws.Range("K120").Value = 10297.5
$ws.Range("L120").Value = 30000
$ws.Range("M120").Value = -5459.5
$ws.Range("N120").Value = -39676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5557778
$ws.Range("I80").Value = 2262.5
$ws.Range("K80").Value = 2262.5
$ws.Range("M80").Value = -1264.5
$ws.Range("H83").Value = 5557778
$ws.Range("I83").Value = 2262.5
$ws.Range("K83").Value = 11312.5
$ws.Range("M83").Value = -6320.5
$ws.Range("H107").Value = 319.8
$ws.Range("I107").Value = 188.18182
$ws.Range("J107").Value = 480.66666
$ws.Range("K107").Value = 188.18182
$ws.Range("L107").Value = 480.66666
$ws.Range("M107").Value = 1731.81818
$ws.Range("N107").Value = -4320.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 25643672
$ws.Range("I55").Value = 5341.6313
$ws.Range("J55").Value = 50000090
$ws.Range("K55").Value = 5341.6313
$ws.Range("L55").Value = 50000090
$ws.Range("M55").Value = -5168.6313
$ws.Range("N55").Value = -50000436
$ws.Range("H68").Value = 1480.8422
$ws.Range("I68").Value = 1426.2354
$ws.Range("J68").Value = 1945
$ws.Range("K68").Value = 1426.2354
$ws.Range("L68").Value = 1945
$ws.Range("M68").Value = -677.2354
$ws.Range("N68").Value = -3443
$ws.Range("H71").Value = 1480.8422
$ws.Range("I71").Value = 1426.2354
$ws.Range("J71").Value = 1945
$ws.Range("K71").Value = 7131.177
$ws.Range("L71").Value = 9725
$ws.Range("M71").Value = -3387.177
$ws.Range("N71").Value = -17213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3180
$ws.Range("I62").Value = 3114.2856
$ws.Range("J62").Value = 3333.3333
$ws.Range("K62").Value = 3114.2856
$ws.Range("L62").Value = 3333.3333
$ws.Range("M62").Value = -2490.2856
$ws.Range("N62").Value = -4581.3333
$ws.Range("H65").Value = 3180
$ws.Range("I65").Value = 3114.2856
$ws.Range("J65").Value = 3333.3333
$ws.Range("K65").Value = 15571.428
$ws.Range("L65").Value = 16666.6665
$ws.Range("M65").Value = -12451.428
$ws.Range("N65").Value = -22906.6665
